# TestItems_STH.xlsx - add a new order-line / API test row.
# A new row (item 121486, qty 1) is inserted right after the existing
# 121486 row (row 9), which also gets its quantities corrected. Every
# row below the insertion point shifts down by one automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the ordered-qty values on the two existing order lines.
$ws.Range("J8").Value = 3
$ws.Range("J9").Value = 2

# Insert a new row for the additional order line (pushes everything
# from row 10 down to row 11, etc.).
$ws.Rows("10:10").Insert()

# Populate the new order line - same item/shipment info as the row
# above it, with quantity 1.
$ws.Range("A10").Value = 121486
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "CnC"
$ws.Range("D10").Value = "PICKUP_IN_STORE"
$ws.Range("E10").Value = "PICK"
$ws.Range("F10").Value = 11990
$ws.Range("G10").Value = 11990
$ws.Range("H10").Value = 11990
$ws.Range("I10").Value = 121486
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 457
$ws.Range("L10").Value = 457

# Match the author's final selection/view state.
$ws.Range("J10").Select() | Out-Null
